# Update the maximum cost driven retirements cap ($B$12 on the "About" sheet)
# from 0.25 to 0.325. All dependent formulas (About!B20:B121, the
# CSC-CSCCCMvSoECBtY row-1 formulas, and the chart caches that reference
# them) will be recalculated automatically by Excel.

$wb = $excel.ActiveWorkbook
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("B12").Value = 0.325

$excel.CalculateFullRebuild()
$wb.RefreshAll()
